$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# D3: "25 minutes" -> "50 minutes"
$ws.Range("D3").Value = "50 minutes"

# C4: 2 -> 3 (this ripples into the SUM total at C31)
$ws.Range("C4").Value = 3

# E3: the note now reads as done/superseded for its first clause, so strike
# through just "Research mobile applications;" and leave the rest normal
$ws.Range("E3").Value = "Research mobile applications; determine PRL-appropriate structure; maybe ask Federica for opinion"
$strikeLen = "Research mobile applications;".Length
$ws.Range("E3").Characters(1, $strikeLen).Font.Strikethrough = $true

# Move the active selection to B3 (matches the recorded view state)
$ws.Range("B3").Select()
